$d = $word.ActiveDocument

$replacements = @(
    @("488÷4=122, 0", "235÷7=33, 4"),
    @("393÷3=131, 0", "361÷6=60, 1"),
    @("459÷4=114, 3", "519÷7=74, 1"),
    @("847÷8=105, 7", "421÷5=84, 1"),
    @("394÷4=98, 2", "575÷7=82, 1"),
    @("288÷8=36, 0", "854÷5=170, 4"),
    @("638÷6=106, 2", "689÷4=172, 1"),
    @("207÷9=23, 0", "684÷9=76, 0"),
    @("958÷2=479, 0", "375÷7=53, 4"),
    @("720÷8=90, 0", "179÷7=25, 4"),
    @("846÷8=105, 6", "474÷3=158, 0"),
    @("204÷5=40, 4", "258÷5=51, 3"),
    @("731÷5=146, 1", "998÷8=124, 6"),
    @("700÷8=87, 4", "170÷7=24, 2"),
    @("152÷3=50, 2", "999÷2=499, 1"),
    @("758÷5=151, 3", "246÷3=82, 0"),
    @("526÷2=263, 0", "952÷6=158, 4"),
    @("686÷6=114, 2", "822÷9=91, 3"),
    @("301÷5=60, 1", "529÷8=66, 1"),
    @("316÷9=35, 1", "986÷9=109, 5"),
    @("439÷8=54, 7", "968÷3=322, 2"),
    @("923÷9=102, 5", "871÷2=435, 1"),
    @("111÷9=12, 3", "803÷9=89, 2"),
    @("404÷4=101, 0", "433÷3=144, 1"),
    @("903÷6=150, 3", "371÷6=61, 5")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
